$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.176273584365845
$ws.Range("B1").Value = 2.413061857223511
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.33787989616394
$ws.Range("E1").Value = 1.201970219612122
